$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 395.08
$ws.Range("I33").Value = 197.10527
$ws.Range("K33").Value = 197.10527
$ws.Range("M33").Value = 31.89473000000001

$ws.Range("H76").Value = 3235.205
$ws.Range("I76").Value = 3235.205
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3235.205
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2920.205
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3235.205
$ws.Range("I79").Value = 3235.205
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3235.205
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2143.205
$ws.Range("N79").ClearContents()

$ws.Range("H101").Value = 10526975
$ws.Range("I101").Value = 18182132
$ws.Range("J101").Value = 1135
$ws.Range("K101").Value = 54546396
$ws.Range("L101").Value = 3405
$ws.Range("M101").Value = -54544774
$ws.Range("N101").Value = -6649

$ws.Range("H111").Value = 1359
$ws.Range("I111").Value = 1396.5
$ws.Range("J111").Value = 1314
$ws.Range("K111").Value = 4189.5
$ws.Range("L111").Value = 3942
$ws.Range("M111").Value = -1122.5
$ws.Range("N111").Value = -10076

$ws.Range("H129").Value = 866.4194
$ws.Range("I129").Value = 557.7778
$ws.Range("J129").Value = 992.6818
$ws.Range("K129").Value = 1673.3334
$ws.Range("L129").Value = 2978.0454
$ws.Range("M129").Value = 3326.6666
$ws.Range("N129").Value = -12978.0454

$ws.Range("H138").Value = 2271.29
$ws.Range("I138").Value = 1214.6364
$ws.Range("J138").Value = 2569.3206
$ws.Range("K138").Value = 3643.9092
$ws.Range("L138").Value = 7707.9618
$ws.Range("M138").Value = 1496.0908
$ws.Range("N138").Value = -17987.9618

$ws.Range("H141").Value = 1786.5454
$ws.Range("I141").Value = 1554.7
$ws.Range("K141").Value = 4664.1
$ws.Range("M141").Value = 515.8999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5335.2
$ws.Range("I32").Value = 4653.718
$ws.Range("J32").Value = 9196.933999999999
$ws.Range("K32").Value = 4653.718
$ws.Range("L32").Value = 9196.933999999999
$ws.Range("M32").Value = -4366.718
$ws.Range("N32").Value = -9770.933999999999

$ws.Range("H45").Value = 993.7143
$ws.Range("J45").Value = 962.4
$ws.Range("L45").Value = 962.4
$ws.Range("N45").Value = -1716.4

$ws.Range("H61").Value = 1848.9592
$ws.Range("I61").Value = 1720.4872
$ws.Range("J61").Value = 2350
$ws.Range("K61").Value = 1720.4872
$ws.Range("L61").Value = 2350
$ws.Range("M61").Value = -1508.4872
$ws.Range("N61").Value = -2774

$ws.Range("H74").Value = 30837.176
$ws.Range("I74").Value = 49021.094
$ws.Range("J74").Value = 1463.1538
$ws.Range("K74").Value = 49021.094
$ws.Range("L74").Value = 1463.1538
$ws.Range("M74").Value = -48147.094
$ws.Range("N74").Value = -3211.1538

$ws.Range("H77").Value = 30837.176
$ws.Range("I77").Value = 49021.094
$ws.Range("J77").Value = 1463.1538
$ws.Range("K77").Value = 245105.47
$ws.Range("L77").Value = 7315.769
$ws.Range("M77").Value = -240737.47
$ws.Range("N77").Value = -16051.769

$ws.Range("H136").Value = 1848.9592
$ws.Range("I136").Value = 1720.4872
$ws.Range("J136").Value = 2350
$ws.Range("K136").Value = 5161.461600000001
$ws.Range("L136").Value = 7050
$ws.Range("M136").Value = -2611.461600000001
$ws.Range("N136").Value = -12150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1459.8518
$ws.Range("I20").Value = 1454.8636
$ws.Range("J20").Value = 1481.8
$ws.Range("K20").Value = 1454.8636
$ws.Range("L20").Value = 1481.8
$ws.Range("M20").Value = -1207.8636
$ws.Range("N20").Value = -1975.8

$ws.Range("H86").Value = 501807.66
$ws.Range("I86").Value = 1866.25
$ws.Range("J86").Value = 1168396.1
$ws.Range("K86").Value = 1866.25
$ws.Range("L86").Value = 1168396.1
$ws.Range("M86").Value = -743.25
$ws.Range("N86").Value = -1170642.1

$ws.Range("H89").Value = 501807.66
$ws.Range("I89").Value = 1866.25
$ws.Range("J89").Value = 1168396.1
$ws.Range("K89").Value = 9331.25
$ws.Range("L89").Value = 5841980.5
$ws.Range("M89").Value = -3715.25
$ws.Range("N89").Value = -5853212.5

$ws.Range("H104").Value = 62684
$ws.Range("J104").Value = 62684
$ws.Range("L104").Value = 62684
$ws.Range("N104").Value = -69672

$ws.Range("H105").Value = 3983005
$ws.Range("I105").Value = 7961010
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 7961010
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -7959263
$ws.Range("N105").Value = -8494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1881.625
$ws.Range("I58").Value = 1421.1428
$ws.Range("J58").Value = 5105
$ws.Range("K58").Value = 1421.1428
$ws.Range("L58").Value = 5105
$ws.Range("M58").Value = -1218.1428
$ws.Range("N58").Value = -5511

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H134").Value = 1312.1346
$ws.Range("I134").Value = 1105.1777
$ws.Range("J134").Value = 2642.5715
$ws.Range("K134").Value = 3315.5331
$ws.Range("L134").Value = 7927.7145
$ws.Range("M134").Value = -780.5330999999996
$ws.Range("N134").Value = -12997.7145

$ws.Range("H136").Value = 1881.625
$ws.Range("I136").Value = 1421.1428
$ws.Range("J136").Value = 5105
$ws.Range("K136").Value = 4263.428400000001
$ws.Range("L136").Value = 15315
$ws.Range("M136").Value = -1713.428400000001
$ws.Range("N136").Value = -20415

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2956.6667
$ws.Range("I125").Value = 495
$ws.Range("J125").Value = 4187.5
$ws.Range("K125").Value = 1485
$ws.Range("L125").Value = 12562.5
$ws.Range("M125").Value = 3435
$ws.Range("N125").Value = -22402.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1119.55
$ws.Range("I97").Value = 1109
$ws.Range("J97").Value = 1130.1
$ws.Range("K97").Value = 1109
$ws.Range("L97").Value = 1130.1
$ws.Range("M97").Value = -613
$ws.Range("N97").Value = -2122.1

$ws.Range("H132").Value = 47560.137
$ws.Range("I132").Value = 1600.7858
$ws.Range("J132").Value = 127989
$ws.Range("K132").Value = 4802.357400000001
$ws.Range("L132").Value = 383967
$ws.Range("M132").Value = -2272.357400000001
$ws.Range("N132").Value = -389027

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3261.7273
$ws.Range("I61").Value = 3590.8
$ws.Range("J61").Value = 2987.5
$ws.Range("K61").Value = 3590.8
$ws.Range("L61").Value = 2987.5
$ws.Range("M61").Value = -3388.8
$ws.Range("N61").Value = -3391.5

$ws.Range("H100").Value = 79212.53999999999
$ws.Range("I100").Value = 168315.5
$ws.Range("J100").Value = 2838.5715
$ws.Range("K100").Value = 168315.5
$ws.Range("L100").Value = 2838.5715
$ws.Range("M100").Value = -167774.5
$ws.Range("N100").Value = -3920.5715

$ws.Range("H113").Value = 3261.7273
$ws.Range("I113").Value = 3590.8
$ws.Range("J113").Value = 2987.5
$ws.Range("K113").Value = 3590.8
$ws.Range("L113").Value = 2987.5
$ws.Range("M113").Value = -1420.8
$ws.Range("N113").Value = -7327.5

$ws.Range("H132").Value = 202403.6
$ws.Range("I132").Value = 53499.69
$ws.Range("J132").Value = 478939.44
$ws.Range("K132").Value = 160499.07
$ws.Range("L132").Value = 1436818.32
$ws.Range("M132").Value = -157969.07
$ws.Range("N132").Value = -1441878.32

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 47724.09
$ws.Range("J46").Value = 47724.09
$ws.Range("L46").Value = 47724.09
$ws.Range("N46").Value = -48186.09

$ws.Range("H132").Value = 4853.107
$ws.Range("I132").Value = 678.3158
$ws.Range("J132").Value = 13666.556
$ws.Range("K132").Value = 2034.9474
$ws.Range("L132").Value = 40999.66800000001
$ws.Range("M132").Value = 495.0526
$ws.Range("N132").Value = -46059.66800000001

$ws.Range("H134").Value = 47724.09
$ws.Range("J134").Value = 47724.09
$ws.Range("L134").Value = 143172.27
$ws.Range("N134").Value = -148242.27
